$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.489.95'
$ws.Range('E2').Value = '  -4.25%  '
$ws.Range('D3').Value = '3.357.01'
$ws.Range('E3').Value = '  -5.06%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '182.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.599'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '3.350.40'
$ws.Range('E9').Value = '  -4.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.186'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.593'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.66'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -7.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000268'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.70'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.05%  '
$ws.Range('D15').Value = '3.897.70'
$ws.Range('E15').Value = '  -4.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '603.34'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -9.25%  '
$ws.Range('D17').Value = '66.601.23'
$ws.Range('E17').Value = '  -4.10%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.363.24'
$ws.Range('E18').Value = '  -4.67%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.02'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.72%  '
$ws.Range('E20').Value = '  -3.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.46'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.911'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '16.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -8.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.12'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.84'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.07'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.10%  '
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.82%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.35'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.74'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -9.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.48'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -8.66%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.28'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.38%  '
$ws.Range('B33').Value = 'dogwifhat'
$ws.Range('C33').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.82'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -14.39%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '561.86'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +11.16%  '
$ws.Range('B35').Value = 'Cosmos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.10'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.58%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '3.835.09'
$ws.Range('E36').Value = '  +0.84%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.105'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '58.13'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.40'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.70%  '
$ws.Range('D41').Value = '0.0₃0716'
$ws.Range('E41').Value = '  -11.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.39'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +22.90%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.127'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.65'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -9.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.346'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '32.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0415'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.15'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.65'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.130'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.12%  '
